# Refined metadata to be additional tab
$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("data")

# --- Update the "time_taken" timestamps on the data sheet (column F, rows 2-18) ---
$ws1.Range("F2").Value  = "2021-10-05 14:22:49.526700"
$ws1.Range("F3").Value  = "2021-10-05 14:22:49.526708"
$ws1.Range("F4").Value  = "2021-10-05 14:22:49.526711"
$ws1.Range("F5").Value  = "2021-10-05 14:22:49.526713"
$ws1.Range("F6").Value  = "2021-10-05 14:22:49.526716"
$ws1.Range("F7").Value  = "2021-10-05 14:22:49.526719"
$ws1.Range("F8").Value  = "2021-10-05 14:22:49.526721"
$ws1.Range("F9").Value  = "2021-10-05 14:22:49.526724"
$ws1.Range("F10").Value = "2021-10-05 14:22:49.526727"
$ws1.Range("F11").Value = "2021-10-05 14:22:49.526729"
$ws1.Range("F12").Value = "2021-10-05 14:22:49.526732"
$ws1.Range("F13").Value = "2021-10-05 14:22:49.526734"
$ws1.Range("F14").Value = "2021-10-05 14:22:49.526737"
$ws1.Range("F15").Value = "2021-10-05 14:22:49.526739"
$ws1.Range("F16").Value = "2021-10-05 14:22:49.526742"
$ws1.Range("F17").Value = "2021-10-05 14:22:49.526744"
$ws1.Range("F18").Value = "2021-10-05 14:22:49.526747"

# --- Add a new "metadata" worksheet right after "data" ---
$wsMeta = $wb.Worksheets.Add($null, $ws1)
$wsMeta.Name = "metadata"

# Header row
$wsMeta.Range("B1").Value = "data_name"
$wsMeta.Range("C1").Value = "data_id"
$wsMeta.Range("D1").Value = "data_version"
$wsMeta.Range("E1").Value = "data_version_created"
$wsMeta.Range("F1").Value = "panel_query_time"
$wsMeta.Range("G1").Value = "panel_get_request"

# Data row
$wsMeta.Range("A2").Value = 0
$wsMeta.Range("B2").Value = "Skeletal muscle channelopathy"
$wsMeta.Range("C2").Value = 542

# data_version must stay a text value ("1.31"), not get coerced to a number
$wsMeta.Range("D2").Value = "'1.31"
$wsMeta.Range("D2").Style = "Normal"

$wsMeta.Range("E2").Value = "2021-07-09T12:42:29.527422Z"
$wsMeta.Range("F2").Value = "2021-10-05 14:22:49.522990"
$wsMeta.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/542/?format=json"

# Match the header style used on the "data" sheet's header row (bold, bordered,
# centered) for the metadata header row and the A2 index cell, by copying the
# existing format instead of synthesizing a brand-new style entry.
$ws1.Range("B1").Copy()
$wsMeta.Range("B1:G1").PasteSpecial(-4122)  # xlPasteFormats
$wsMeta.Range("A2").PasteSpecial(-4122)     # xlPasteFormats

Write-Host "metadata sheet added"
